$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("秋元")
$r.Text = "秋本"
